$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 2450.2
$ws.Range("J7").Value = 1062.75
$ws.Range("L7").Value = 1062.75
$ws.Range("N7").Value = -1286.75
$ws.Range("H10").Value = 27.5
$ws.Range("I10").Value = 27.5
$ws.Range("K10").Value = 27.5
$ws.Range("M10").Value = 265.5
$ws.Range("H14").Value = 2450.2
$ws.Range("J14").Value = 1062.75
$ws.Range("L14").Value = 1062.75
$ws.Range("N14").Value = -1444.75
$ws.Range("H20").Value = 1599.1666
$ws.Range("I20").Value = 1319
$ws.Range("J20").Value = 3000
$ws.Range("K20").Value = 1319
$ws.Range("L20").Value = 3000
$ws.Range("M20").Value = -1089
$ws.Range("N20").Value = -3460
$ws.Range("H32").Value = 3148.5
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 4297
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 4297
$ws.Range("M32").Value = -1674
$ws.Range("N32").Value = -4949
$ws.Range("H35").Value = 1599.1666
$ws.Range("I35").Value = 1319
$ws.Range("J35").Value = 3000
$ws.Range("K35").Value = 1319
$ws.Range("L35").Value = 3000
$ws.Range("M35").Value = -940
$ws.Range("N35").Value = -3758
$ws.Range("H39").Value = 288.08334
$ws.Range("I39").Value = 286.9
$ws.Range("J39").Value = 294
$ws.Range("K39").Value = 860.6999999999999
$ws.Range("L39").Value = 882
$ws.Range("M39").Value = -564.6999999999999
$ws.Range("N39").Value = -1474
$ws.Range("H41").Value = 2029.3334
$ws.Range("I41").Value = 1401
$ws.Range("J41").Value = 2086.4546
$ws.Range("K41").Value = 1401
$ws.Range("L41").Value = 2086.4546
$ws.Range("M41").Value = -961
$ws.Range("N41").Value = -2966.4546
$ws.Range("H74").Value = 9915
$ws.Range("I74").Value = 7898
$ws.Range("K74").Value = 7898
$ws.Range("M74").Value = -6962
$ws.Range("H77").Value = 9915
$ws.Range("I77").Value = 7898
$ws.Range("K77").Value = 39490
$ws.Range("M77").Value = -34810
$ws.Range("H87").Value = 77554
$ws.Range("J87").Value = 92943
$ws.Range("L87").Value = 92943
$ws.Range("N87").Value = -95439
$ws.Range("H90").Value = 77554
$ws.Range("J90").Value = 92943
$ws.Range("L90").Value = 278829
$ws.Range("N90").Value = -291309
$ws.Range("H132").Value = 3358.75
$ws.Range("I132").Value = 3146
$ws.Range("K132").Value = 9438
$ws.Range("M132").Value = -6908
$ws.Range("H137").Value = 2685
$ws.Range("I137").Value = 1992.1428
$ws.Range("K137").Value = 5976.428400000001
$ws.Range("M137").Value = -3426.428400000001

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 822.76
$ws.Range("I2").Value = 711.3333
$ws.Range("J2").Value = 989.9
$ws.Range("K2").Value = 711.3333
$ws.Range("L2").Value = 989.9
$ws.Range("M2").Value = -598.3333
$ws.Range("N2").Value = -1215.9
$ws.Range("H45").Value = 3097.9285
$ws.Range("I45").Value = 2232.8333
$ws.Range("J45").Value = 3746.75
$ws.Range("K45").Value = 2232.8333
$ws.Range("L45").Value = 3746.75
$ws.Range("M45").Value = -1855.8333
$ws.Range("N45").Value = -4500.75
$ws.Range("H61").Value = 3559.8
$ws.Range("I61").Value = 3425.9167
$ws.Range("J61").Value = 4095.3333
$ws.Range("K61").Value = 3425.9167
$ws.Range("L61").Value = 4095.3333
$ws.Range("M61").Value = -3213.9167
$ws.Range("N61").Value = -4519.3333
$ws.Range("H116").Value = 822.76
$ws.Range("I116").Value = 711.3333
$ws.Range("J116").Value = 989.9
$ws.Range("K116").Value = 711.3333
$ws.Range("L116").Value = 989.9
$ws.Range("M116").Value = 1582.6667
$ws.Range("N116").Value = -5577.9
$ws.Range("H132").Value = 2386.6
$ws.Range("I132").Value = 2006.5
$ws.Range("K132").Value = 6019.5
$ws.Range("M132").Value = -3489.5
$ws.Range("H136").Value = 3559.8
$ws.Range("I136").Value = 3425.9167
$ws.Range("J136").Value = 4095.3333
$ws.Range("K136").Value = 10277.7501
$ws.Range("L136").Value = 12285.9999
$ws.Range("M136").Value = -7727.750100000001
$ws.Range("N136").Value = -17385.9999
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 822.76
$ws.Range("I3").Value = 711.3333
$ws.Range("J3").Value = 989.9
$ws.Range("K3").Value = 711.3333
$ws.Range("L3").Value = 989.9
$ws.Range("M3").Value = -597.3333
$ws.Range("N3").Value = -1217.9
$ws.Range("H20").Value = 1089.8462
$ws.Range("I20").Value = 933
$ws.Range("J20").Value = 1952.5
$ws.Range("K20").Value = 933
$ws.Range("L20").Value = 1952.5
$ws.Range("M20").Value = -686
$ws.Range("N20").Value = -2446.5
$ws.Range("H80").Value = 693.4167
$ws.Range("I80").Value = 437.4
$ws.Range("K80").Value = 437.4
$ws.Range("M80").Value = 560.6
$ws.Range("H83").Value = 693.4167
$ws.Range("I83").Value = 437.4
$ws.Range("K83").Value = 2187
$ws.Range("M83").Value = 2805
$ws.Range("H94").Value = 1121.3182
$ws.Range("I94").Value = 1120.0526
$ws.Range("K94").Value = 1120.0526
$ws.Range("M94").Value = -669.0526
$ws.Range("H107").Value = 3845.6924
$ws.Range("I107").Value = 2080.8333
$ws.Range("K107").Value = 2080.8333
$ws.Range("M107").Value = -160.8332999999998

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 5000
$ws.Range("J36").Value = 5000
$ws.Range("L36").Value = 5000
$ws.Range("N36").Value = -5776
$ws.Range("H40").Value = 5000
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5320
$ws.Range("H58").Value = 2499
$ws.Range("I58").Value = 1439.6
$ws.Range("K58").Value = 1439.6
$ws.Range("M58").Value = -1236.6
$ws.Range("H135").Value = 78015
$ws.Range("I135").Value = 100500
$ws.Range("J135").Value = 70520
$ws.Range("K135").Value = 100500
$ws.Range("L135").Value = 70520
$ws.Range("M135").Value = -95430
$ws.Range("N135").Value = -80660
$ws.Range("H136").Value = 2499
$ws.Range("I136").Value = 1439.6
$ws.Range("K136").Value = 4318.799999999999
$ws.Range("M136").Value = -1768.799999999999

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 707449.3
$ws.Range("I4").Value = 3002247.8
$ws.Range("J4").Value = 1357.4615
$ws.Range("K4").Value = 9006743.399999999
$ws.Range("L4").Value = 4072.3845
$ws.Range("M4").Value = -9006631.399999999
$ws.Range("N4").Value = -4296.3845
$ws.Range("H7").Value = 306.84616
$ws.Range("I7").Value = 293.16666
$ws.Range("K7").Value = 879.4999799999999
$ws.Range("M7").Value = -767.4999799999999
$ws.Range("H9").Value = 317.8
$ws.Range("J9").Value = 761.5
$ws.Range("L9").Value = 2284.5
$ws.Range("N9").Value = -2732.5
$ws.Range("H14").Value = 457.07144
$ws.Range("I14").Value = 457.07144
$ws.Range("K14").Value = 1371.21432
$ws.Range("M14").Value = -1198.21432
$ws.Range("H34").Value = 726.9091
$ws.Range("J34").Value = 1035.5714
$ws.Range("L34").Value = 3106.7142
$ws.Range("N34").Value = -3274.7142
$ws.Range("H44").Value = 10000
$ws.Range("J44").Value = 10000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30796
$ws.Range("H50").Value = 1484.75
$ws.Range("I50").Value = 713
$ws.Range("J50").Value = 3800
$ws.Range("K50").Value = 2139
$ws.Range("L50").Value = 11400
$ws.Range("M50").Value = -1658
$ws.Range("N50").Value = -12362
$ws.Range("H53").Value = 1484.75
$ws.Range("I53").Value = 713
$ws.Range("J53").Value = 3800
$ws.Range("K53").Value = 2139
$ws.Range("L53").Value = 11400
$ws.Range("M53").Value = -1658
$ws.Range("N53").Value = -12362
$ws.Range("H132").Value = 1800
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -15470
$ws.Range("N132").Value = -18560

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3499.5
$ws.Range("I70").Value = 2999
$ws.Range("K70").Value = 2999
$ws.Range("M70").Value = -2729
$ws.Range("H73").Value = 3499.5
$ws.Range("I73").Value = 2999
$ws.Range("K73").Value = 2999
$ws.Range("M73").Value = -2063
$ws.Range("H102").Value = 1148
$ws.Range("I102").Value = 1191.2667
$ws.Range("K102").Value = 1191.2667
$ws.Range("M102").Value = 430.7333000000001
$ws.Range("H132").Value = 5254.75
$ws.Range("I132").Value = 6200
$ws.Range("K132").Value = 18600
$ws.Range("M132").Value = -16070

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = $null
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = $null
$ws.Range("N126").Value = $null
$ws.Range("H132").Value = 5231.636
$ws.Range("I132").Value = 4931.5
$ws.Range("K132").Value = 14794.5
$ws.Range("M132").Value = -12264.5

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 84500
$ws.Range("J46").Value = 84500
$ws.Range("L46").Value = 84500
$ws.Range("N46").Value = -84962
$ws.Range("H134").Value = 84500
$ws.Range("J134").Value = 84500
$ws.Range("L134").Value = 253500
$ws.Range("N134").Value = -258570
$ws.Range("H136").Value = 3283.3845
$ws.Range("I136").Value = 2421.2222
$ws.Range("K136").Value = 7263.6666
$ws.Range("M136").Value = -4713.6666
